$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.821.20"
$ws.Range("E2").Value = "  -1.45%  "
$ws.Range("D3").Value = "1.873.53"
$ws.Range("E3").Value = "  -1.63%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "'301.14"
$ws.Range("D5").ClearFormats()
$ws.Range("D6").Value = "'1.001"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.10%  "
$ws.Range("D7").Value = "'0.5339"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +1.67%  "
$ws.Range("D8").Value = "'0.3761"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -1.34%  "
$ws.Range("D9").Value = "'0.07181"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -1.70%  "
$ws.Range("D10").Value = "'21.62"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.46%  "
$ws.Range("D11").Value = "'0.8875"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -2.02%  "
$ws.Range("D12").Value = "'0.08174"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +1.35%  "
$ws.Range("D13").Value = "'93.56"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -2.29%  "
$ws.Range("D14").Value = "1.826.86"
$ws.Range("E14").Value = "  -1.77%  "
$ws.Range("D15").Value = "'5.283"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -1.42%  "
$ws.Range("E16").Value = "  -0.14%  "
$ws.Range("D17").Value = "'14.78"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.34%  "
$ws.Range("D18").Value = "'0.000008547"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -1.59%  "
$ws.Range("D19").Value = "'1.001"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.10%  "
$ws.Range("D20").Value = "26.838.19"
$ws.Range("E20").Value = "  -1.52%  "
$ws.Range("E21").Value = "  -2.78%  "
$ws.Range("E22").Value = "  -1.32%  "
$ws.Range("D23").Value = "'6.394"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -1.30%  "
$ws.Range("D24").Value = "'146.38"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -1.94%  "
$ws.Range("D25").Value = "'2.267"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -3.98%  "
$ws.Range("D26").Value = "'1.735"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.42%  "
$ws.Range("D27").Value = "'18.02"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -1.27%  "
$ws.Range("D28").Value = "'113.82"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -2.36%  "
$ws.Range("D29").Value = "'4.727"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -2.35%  "
$ws.Range("E30").Value = "  -5.98%  "
$ws.Range("D31").Value = "'0.09146"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.97%  "
$ws.Range("D32").Value = "'0.8069"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.69%  "
$ws.Range("E33").Value = "  -2.25%  "
$ws.Range("E34").Value = "  -4.58%  "
$ws.Range("D35").Value = "'2.966"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.04%  "
$ws.Range("D36").Value = "'0.5999"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +4.70%  "
$ws.Range("D37").Value = "'3.194"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -5.33%  "
$ws.Range("D38").Value = "'2.587"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -3.53%  "
$ws.Range("D39").Value = "'0.01950"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -2.13%  "
$ws.Range("D40").Value = "'1.072"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -1.43%  "
$ws.Range("D41").Value = "'8.871"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -1.60%  "
$ws.Range("D42").Value = "'6.554"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.67%  "
$ws.Range("B43").Value = "Decentraland"
$ws.Range("C43").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D43").Value = "'0.5129"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +4.64%  "
$ws.Range("B44").Value = "Quant"
$ws.Range("C44").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D44").Value = "'114.75"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -1.73%  "
$ws.Range("E45").Value = "  -1.94%  "
$ws.Range("E46").Value = "  -0.14%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "'9.948"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -2.02%  "
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").Value = "'1.631"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.31%  "
$ws.Range("D49").Value = "'37.52"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -2.79%  "
$ws.Range("D50").Value = "'0.06050"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +1.50%  "
$ws.Range("D51").Value = "'62.05"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -3.66%  "
